{"js": "// Locate the paragraphs we need to edit by their current text content so the\n// script is resilient to exact indices. We still rely on the paragraph\n// collection (body.paragraphs) which mirrors the document order.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet fromPara = null;\nlet designPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"From midpoint between the two origins\") {\n    fromPara = paragraphs.items[i];\n  } else if (t === \"Design\") {\n    designPara = paragraphs.items[i];\n  }\n}\n\nif (!fromPara) {\n  throw new Error(\"Could not find the 'From ... midpoint between the two origins' paragraph\");\n}\n\n// Merge the two runs (\"From \" + \"midpoint between the two origins\") into a\n// single run with the new, expanded search-criteria text. insertText with\n// \"Replace\" rewrites the whole paragraph range's text as one run while\n// keeping the paragraph's own formatting (list numbering, style, etc.).\nfromPara.insertText(\n  \"Distance: From midpoint between the two origins, select x km radius\",\n  \"Replace\"\n);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark from the end of \"Purple: #5D5390\" to the start\n// of the \"Design\" paragraph. Remove the old one first, then insert the new\n// one, so both mutations target the same bookmark name without colliding.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nif (designPara) {\n  const designStart = designPara.getRange(\"Start\");\n  designStart.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Locate the paragraphs we care about by their text -------------------\n$fromPara = $null\n$designPara = $null\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -eq \"From midpoint between the two origins`r\") {\n        $fromPara = $p\n    } elseif ($t -eq \"Design`r\") {\n        $designPara = $p\n    }\n}\n\nif ($fromPara -eq $null) {\n    throw \"Could not find the 'From ... midpoint between the two origins' paragraph\"\n}\n\n# --- Merge the two runs (\"From \" + \"midpoint between the two origins\") ---\n# into a single run with the new, expanded search-criteria text. Exclude the\n# trailing paragraph mark from the replaced range so the paragraph itself\n# (and its list/number formatting) is preserved, only its run content changes.\n$start = $fromPara.Range.Start\n$end = $fromPara.Range.End - 1\n$textRange = $d.Range($start, $end)\n$textRange.Text = \"Distance: From midpoint between the two origins, select x km radius\"\n\n# --- Move the \"_GoBack\" bookmark from the end of \"Purple: #5D5390\" to the --\n# start of the \"Design\" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\nif ($designPara -ne $null) {\n    $designStart = $d.Range($designPara.Range.Start, $designPara.Range.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $designStart)\n}\n"}
